$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("resumo")

$ws.Rows.Item(26).Insert()
$ws.Rows.Item(27).Insert()
$ws.Rows.Item(28).Insert()
$ws.Rows.Item(29).Insert()

$ws.Range("H31:I31").ClearContents()
$ws.Range("E31").NumberFormat = "0%"

$data = @(
  @(43929, 1376, 53),
  @(43930, 1445, 57),
  @(43931, 1558, 67),
  @(43932, 1668, 74)
)
$r = 26
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Formula = "=C$r/B$r"
  $ws.Cells.Item($r, 5).Formula = "=B$r/B$($r-1) - 1"
  $ws.Cells.Item($r, 6).Formula = "=10^6*B$r/`$Q`$1"
  $ws.Cells.Item($r, 7).Formula = "=B$r-B$($r-1)"
  $ws.Cells.Item($r, 9).Formula = "=B$r/B$($r-1)"
  $r++
}

$ws.Activate()
$ws.Range("D12").Select()
